# Aakash CLUB REMOVED AS PER PARTH SIR EMAIL
#
# The "Aakash" club occupied row 3 of Sheet1 (Club ID "Aakash", Club Name
# "Adwitya Aakash Club", logo "images1\Student Club\ADWITIYA AAKASH - circle
# - dark.webp", ...). Per the author's request, remove that whole row so the
# remaining club rows shift up one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The worksheet's external hyperlinks (club logo links + the H4 gallery
# link) are anchored to fixed cells. Re-create them (minus the one that
# lived on the row being removed) against their post-delete cell positions,
# in the same relative order they were originally defined, so the
# relationship ids come out sequential/minimal just like a manual Excel
# edit would produce.
$hyperlinkTargets = @(
    @{ Cell = "H4"; Url = "images1\gallery\AWS Commuinity Day\aws1.webp" },
    @{ Cell = "C2"; Url = "images1\Student Club\aws.webp" },
    @{ Cell = "C3"; Url = "images1\Student Club\ADWITIYA AAKASH - circle - dark.webp" },
    @{ Cell = "C9"; Url = "images1\Student Club\ML_Club.webp" },
    @{ Cell = "C8"; Url = "images1\Student Club\DS_CLUB.webp" },
    @{ Cell = "C7"; Url = "images1\Student Club\Cyber_Security.webp" },
    @{ Cell = "C4"; Url = "images1\Student Club\Aiml.webp" },
    @{ Cell = "C5"; Url = "images1\Student Club\CC.webp" },
    @{ Cell = "C6"; Url = "images1\Student Club\Club_Gamma.webp" },
    @{ Cell = "C10"; Url = "images1\Student Club\innovation.webp" },
    @{ Cell = "C11"; Url = "images1\Student Club\ChipXplorers_Club.webp" },
    @{ Cell = "C12"; Url = "images1\Student Club\CPSQAUD.webp" },
    @{ Cell = "C13"; Url = "images1\Student Club\eloquence_logo.webp" },
    @{ Cell = "C14"; Url = "images1\Student Club\Scrapping_The_sky.webp" },
    @{ Cell = "C15"; Url = "images1\Student Club\Debate_Club.webp" },
    @{ Cell = "C16"; Url = "images1\Student Club\Code_for_the_cause.webp" }
)

$removedRow = 3

# Clear all existing hyperlinks; they'll be rebuilt below at their
# shifted locations once the row is gone.
$ws.Hyperlinks.Delete()

# Delete the Aakash club's row; everything below shifts up one row.
$ws.Rows("$removedRow`:$removedRow").Delete()

foreach ($link in $hyperlinkTargets) {
    $cell = $link.Cell
    if ($cell -match '^([A-Z]+)(\d+)$') {
        $col = $Matches[1]
        $row = [int]$Matches[2]
    }

    if ($row -eq $removedRow) {
        # This hyperlink lived on the deleted row - drop it entirely.
        continue
    } elseif ($row -gt $removedRow) {
        $row = $row - 1
    }

    $ws.Hyperlinks.Add($ws.Range("$col$row"), $link.Url)
}
